$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Test AUC" values for the two already-existing tuned rows ---
$ws.Cells.Item(52, 15).Value = 0.60299999999999998
$ws.Cells.Item(53, 15).Value = 0.63400000000000001

# --- Grow the table by 7 rows (54-60) ---
$lo = $ws.ListObjects.Item(1)
for ($i = 0; $i -lt 7; $i++) {
    $lo.ListRows.Add() | Out-Null
}

function Set-Row($r, $B, $C, $D, $E, $F, $I, $J, $K, $L, $M) {
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 4).Value = $D
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $L
    $ws.Cells.Item($r, 13).Value = $M
}

Set-Row 54 "RandomForest" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels" 0.63400000000000001 0.75600000000000001 0.65700000000000003 "0.682 (0.053)"
Set-Row 55 "Logistic Classifier" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels" 0.60899999999999999 0.69699999999999995 0.63300000000000001 "0.646 (0.037)"
Set-Row 56 "RandomForest" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels, standard scaling" 0.61199999999999999 0.70699999999999996 0.65200000000000002 "0.657 (0.039)"
Set-Row 57 "RandomForest" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels, minmax scaling" 0.6 0.74299999999999999 0.66100000000000003 "0.668 (0.059)"
Set-Row 58 "XGB Classifier" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels" 0.65400000000000003 0.69599999999999995 0.59899999999999998 "0.649 (0.040)"
Set-Row 59 "XGB Classifier" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels, standard scaling" 0.63 0.71499999999999997 0.59099999999999997 "0.645 (0.052)"
Set-Row 60 "XGB Classifier" "MoCo" "Centers" "1 x 3" "average" "incl. Zoom levels, minmax scaling" 0.64400000000000002 0.63 0.59099999999999997 "0.621 (0.022)"

# --- Extend the conditional-formatting ranges that tracked the table body (not auto-growing with it) ---
$cfs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $cfs.Count; $i++) {
    $fc = $cfs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq '$K$2:$K$53') {
        $fc.ModifyAppliesToRange($ws.Range("K2:K60"))
    }
    if ($addr -eq '$J$2:$J$53') {
        $fc.ModifyAppliesToRange($ws.Range("J2:J60"))
    }
}

# --- Update the view: scrolled position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("O52").Select()
